# "updated DNS resolution diagrams"
#
# Changes applied to the single slide in this deck:
#   1. "www.apps." -> "www." in the "Snip Single Corner of Rectangle 16" shape
#      (keeps the other two runs - "gslb" and ".mycompany.com" - untouched).
#   2. "Cluster's router IP address" -> "Router public IP address" in the
#      "Snip Single Corner of Rectangle 22" shape.
#   3. The two yellow "Right Arrow" callout shapes ("Right Arrow 14" and
#      "Right Arrow 15") are moved from just after the picture to the very
#      end of the z-order (after "TextBox 27" / "F5 Cloud Services DNS LB").
#
# (Datetime "30/04/2020" -> "01/05/2020" footer field edits from the source
# diff belong to other slides of the original deck that aren't present here;
# this single slide has no datetime fields to update.)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# 1) "www.apps." -> "www."
$shApps = $s.Shapes.Item("Snip Single Corner of Rectangle 16")
$runApps = $shApps.TextFrame.TextRange.Runs(1)
$runApps.Text = "www."

# 2) "Cluster's router IP address" -> "Router public IP address"
$shRouter = $s.Shapes.Item("Snip Single Corner of Rectangle 22")
$shRouter.TextFrame.TextRange.Text = "Router public IP address"

# 3) Move the two "Right Arrow" shapes to the end of the shape stack.
$s.Shapes.Item("Right Arrow 14").ZOrder(0)  # msoBringToFront
$s.Shapes.Item("Right Arrow 15").ZOrder(0)  # msoBringToFront
